$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as literal text even when it
# looks numeric (e.g. "0.76" or "641,530,687,575"), without leaving the
# cell's style pointing at a non-default format.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 16
$ws.Range("D16").Value = "image_20250807110238_ppp0.jpg"
Set-TextValue "I16" "641,530,687,575"
Set-TextValue "J16" "0.76"

# Row 17
$ws.Range("D17").Value = "image_20250807110238_ppp0.jpg"
Set-TextValue "I17" "793,481,831,527"
Set-TextValue "J17" "0.71"
